$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.127973556518555
$ws.Range("B1").Value = 2.882835388183594
$ws.Range("C1").Value = 2.224223136901855
$ws.Range("D1").Value = 2.079066276550293
$ws.Range("E1").Value = 2.085942268371582
